$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.318.92'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '2.599.16'
$ws.Range('E3').Value = '  +4.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.58%  '
$ws.Range('E7').Value = '  +4.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +8.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.49%  '
$ws.Range('E11').Value = '  +6.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.17'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.14'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.22%  '
$ws.Range('D14').Value = '2.992.63'
$ws.Range('E14').Value = '  +4.52%  '
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').Value = '2.603.33'
$ws.Range('E16').Value = '  +4.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.918'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.38%  '
$ws.Range('E18').Value = '  +3.70%  '
$ws.Range('D19').Value = '46.383.37'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('E20').Value = '  +6.17%  '
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.73'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '274.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +10.90%  '
$ws.Range('E25').Value = '  +7.07%  '
$ws.Range('E26').Value = '  +9.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +30.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.60'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.23%  '
$ws.Range('E32').Value = '  -1.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.36'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.86'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0839'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '150.73'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.67%  '
$ws.Range('E39').Value = '  +5.54%  '
$ws.Range('E40').Value = '  +5.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.30'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +40.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.84'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0331'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.68%  '
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').Value = '2.136.04'
$ws.Range('E46').Value = '  +5.85%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '92.76'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.78'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.32%  '
